$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 864.7857
$ws.Range("J38").Value = 3419
$ws.Range("L38").Value = 10257
$ws.Range("N38").Value = -11001
$ws.Range("H100").Value = 4496
$ws.Range("I100").Value = 4499.75
$ws.Range("J100").Value = 4488.5
$ws.Range("K100").Value = 4499.75
$ws.Range("L100").Value = 4488.5
$ws.Range("M100").Value = -3958.75
$ws.Range("N100").Value = -5570.5
$ws.Range("H123").Value = 169990
$ws.Range("J123").Value = 169990
$ws.Range("L123").Value = 169990
$ws.Range("N123").Value = -179790
$ws.Range("H125").Value = 1699
$ws.Range("I125").Value = 500
$ws.Range("K125").Value = 4500
$ws.Range("M125").Value = -2040
$ws.Range("H135").Value = 569.1667
$ws.Range("I135").Value = 532.05884
$ws.Range("K135").Value = 4788.52956
$ws.Range("M135").Value = -2253.52956
$ws.Range("H137").Value = 1791548.5
$ws.Range("J137").Value = 4006191
$ws.Range("L137").Value = 12018573
$ws.Range("N137").Value = -12023673

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3008.1343
$ws.Range("I32").Value = 1873.7593
$ws.Range("J32").Value = 7720.154
$ws.Range("K32").Value = 1873.7593
$ws.Range("L32").Value = 7720.154
$ws.Range("M32").Value = -1586.7593
$ws.Range("N32").Value = -8294.154
$ws.Range("H61").Value = 2142210.5
$ws.Range("I61").Value = 83574.84
$ws.Range("J61").Value = 3926361.5
$ws.Range("K61").Value = 83574.84
$ws.Range("L61").Value = 3926361.5
$ws.Range("M61").Value = -83362.84
$ws.Range("N61").Value = -3926785.5
$ws.Range("H74").Value = 460037.97
$ws.Range("I74").Value = 2457.257
$ws.Range("K74").Value = 2457.257
$ws.Range("M74").Value = -1583.257
$ws.Range("H77").Value = 460037.97
$ws.Range("I77").Value = 2457.257
$ws.Range("K77").Value = 12286.285
$ws.Range("M77").Value = -7918.285
$ws.Range("H97").Value = 10813.25
$ws.Range("I97").Value = 10976
$ws.Range("J97").Value = 9999.5
$ws.Range("K97").Value = 10976
$ws.Range("L97").Value = 9999.5
$ws.Range("M97").Value = -10480
$ws.Range("N97").Value = -10991.5
$ws.Range("H103").Value = 55555
$ws.Range("J103").Value = 55555
$ws.Range("L103").Value = 55555
$ws.Range("N103").Value = -57899
$ws.Range("H132").Value = 3125.6843
$ws.Range("I132").Value = 2409
$ws.Range("K132").Value = 7227
$ws.Range("M132").Value = -4697
$ws.Range("H136").Value = 2142210.5
$ws.Range("I136").Value = 83574.84
$ws.Range("J136").Value = 3926361.5
$ws.Range("K136").Value = 250724.52
$ws.Range("L136").Value = 11779084.5
$ws.Range("M136").Value = -248174.52
$ws.Range("N136").Value = -11784184.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H62").Value = 23500
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 23500
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()
$ws.Range("H95").Value = 47795.855
$ws.Range("J95").Value = 49167.332
$ws.Range("L95").Value = 49167.332
$ws.Range("N95").Value = -54659.332
$ws.Range("H105").Value = 9152.939
$ws.Range("I105").Value = 7898.8696
$ws.Range("K105").Value = 7898.8696
$ws.Range("M105").Value = -6151.8696
$ws.Range("H130").Value = 38000
$ws.Range("J130").Value = 38000
$ws.Range("L130").Value = 38000
$ws.Range("N130").Value = -48040
$ws.Range("H134").Value = 90004410
$ws.Range("I134").Value = 4882.75
$ws.Range("K134").Value = 14648.25
$ws.Range("M134").Value = -12113.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2583.1343
$ws.Range("I31").Value = 3801.6924
$ws.Range("J31").Value = 2289.7778
$ws.Range("K31").Value = 3801.6924
$ws.Range("L31").Value = 2289.7778
$ws.Range("M31").Value = -3506.6924
$ws.Range("N31").Value = -2879.7778
$ws.Range("H34").Value = 2583.1343
$ws.Range("I34").Value = 3801.6924
$ws.Range("J34").Value = 2289.7778
$ws.Range("K34").Value = 3801.6924
$ws.Range("L34").Value = 2289.7778
$ws.Range("M34").Value = -3599.6924
$ws.Range("N34").Value = -2693.7778
$ws.Range("H109").Value = 55766.11
$ws.Range("I109").Value = 63979
$ws.Range("K109").Value = 63979
$ws.Range("M109").Value = -62939
$ws.Range("H132").Value = 18522450
$ws.Range("I132").Value = 3665.1667
$ws.Range("K132").Value = 10995.5001
$ws.Range("M132").Value = -8465.500100000001
$ws.Range("H134").Value = 3029.9119
$ws.Range("I134").Value = 2690.7585
$ws.Range("J134").Value = 4997
$ws.Range("K134").Value = 8072.2755
$ws.Range("L134").Value = 14991
$ws.Range("M134").Value = -5537.2755
$ws.Range("N134").Value = -20061

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 657.4808
$ws.Range("J2").Value = 937.4286
$ws.Range("L2").Value = 5624.571599999999
$ws.Range("N2").Value = -5850.571599999999
$ws.Range("H9").Value = 8420018
$ws.Range("I9").Value = 370.16666
$ws.Range("J9").Value = 16839666
$ws.Range("K9").Value = 1110.49998
$ws.Range("L9").Value = 50518998
$ws.Range("M9").Value = -886.4999800000001
$ws.Range("N9").Value = -50519446
$ws.Range("H10").Value = 995.64703
$ws.Range("I10").Value = 565.5
$ws.Range("J10").Value = 1128
$ws.Range("K10").Value = 1696.5
$ws.Range("L10").Value = 3384
$ws.Range("M10").Value = -1557.5
$ws.Range("N10").Value = -3662
$ws.Range("H12").Value = 907.75
$ws.Range("I12").Value = 225
$ws.Range("J12").Value = 1044.3
$ws.Range("K12").Value = 675
$ws.Range("L12").Value = 3132.9
$ws.Range("M12").Value = -502
$ws.Range("N12").Value = -3478.9
$ws.Range("H14").Value = 3203.6875
$ws.Range("I14").Value = 3203.6875
$ws.Range("K14").Value = 9611.0625
$ws.Range("M14").Value = -9438.0625
$ws.Range("H68").Value = 1281.4
$ws.Range("J68").Value = 1175.1364
$ws.Range("L68").Value = 3525.4092
$ws.Range("N68").Value = -5147.4092
$ws.Range("H71").Value = 1281.4
$ws.Range("J71").Value = 1175.1364
$ws.Range("L71").Value = 10576.2276
$ws.Range("N71").Value = -18688.2276
$ws.Range("H107").Value = 2073.8333
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 2073.8333
$ws.Range("K107").Value = 0
$ws.Range("L107").Value = 6221.499899999999
$ws.Range("N107").Value = -10061.4999
$ws.Range("H131").Value = 5613595
$ws.Range("J131").Value = 5557974
$ws.Range("L131").Value = 16673922
$ws.Range("N131").Value = -16684002
$ws.Range("H137").Value = 1568.5714
$ws.Range("I137").Value = 1163.3334
$ws.Range("K137").Value = 3490.0002
$ws.Range("M137").Value = 1609.9998

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 135.14285
$ws.Range("J2").Value = 236.4
$ws.Range("L2").Value = 236.4
$ws.Range("N2").Value = -462.4
$ws.Range("H107").Value = 335166
$ws.Range("I107").Value = 500249
$ws.Range("K107").Value = 500249
$ws.Range("M107").Value = -498329
$ws.Range("H126").Value = 2875.9092
$ws.Range("I126").Value = 2348.6365
$ws.Range("J126").Value = 3403.182
$ws.Range("K126").Value = 7045.9095
$ws.Range("L126").Value = 10209.546
$ws.Range("M126").Value = -4575.9095
$ws.Range("N126").Value = -15149.546
$ws.Range("H132").Value = 11029773
$ws.Range("I132").Value = 3871.4
$ws.Range("J132").Value = 18380374
$ws.Range("K132").Value = 11614.2
$ws.Range("L132").Value = 55141122
$ws.Range("M132").Value = -9084.200000000001
$ws.Range("N132").Value = -55146182
$ws.Range("H134").Value = 77601.875
$ws.Range("J134").Value = 77601.875
$ws.Range("L134").Value = 232805.625
$ws.Range("N134").Value = -237875.625
$ws.Range("H136").Value = 90758.78
$ws.Range("J136").Value = 90758.78
$ws.Range("L136").Value = 272276.34
$ws.Range("N136").Value = -277376.34

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 755.96
$ws.Range("I55").Value = 791.7
$ws.Range("J55").Value = 732.13336
$ws.Range("K55").Value = 791.7
$ws.Range("L55").Value = 732.13336
$ws.Range("M55").Value = -618.7
$ws.Range("N55").Value = -1078.13336
$ws.Range("H132").Value = 5347.857
$ws.Range("I132").Value = 3465.4614
$ws.Range("K132").Value = 10396.3842
$ws.Range("M132").Value = -7866.3842
$ws.Range("H136").Value = 2402.3901
$ws.Range("I136").Value = 2456.8147
$ws.Range("J136").Value = 2297.4285
$ws.Range("K136").Value = 7370.4441
$ws.Range("L136").Value = 6892.2855
$ws.Range("M136").Value = -4820.4441
$ws.Range("N136").Value = -11992.2855

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2076.68
$ws.Range("I132").Value = 1801.5555
$ws.Range("K132").Value = 5404.666499999999
$ws.Range("M132").Value = -2874.666499999999
$ws.Range("H136").Value = 2198.4
$ws.Range("I136").Value = 1222.4166
$ws.Range("K136").Value = 3667.2498
$ws.Range("M136").Value = -1117.2498
